$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the "Added Wumpus..." note with the new sentences about finishing
# the GUI art for pits, breezes, gold/glitter and the random-seed TODO.
$ws.Range("D9").Value = "Added Wumpus and stenches in addition to Adventurer. Renders well, simplify the iconography so I can just use core drawing functions. Solid ovals for things, empty ovals for their signals. Must update documention. Got the rest of the images drawn for the pits, breezes, gold and glitter. Need to turn of the random seed as we get no variation in runs."

# Thursday 21st logged 4.5 hours instead of 1 (F2's SUM(C:C) recalculates automatically).
$ws.Range("C9").Value = 4.5

# The longer note now wraps to more lines, so the row grows taller.
$ws.Rows("9").RowHeight = 57

# Current selection moved from D10 to C10.
$ws.Range("C10").Select() | Out-Null
